# Review.xlsx update per "Update review sheet according to CYRS updates and SRS"
#
# Summary of the edit:
#  - Review sheet: several review points get closed/accepted, a couple of new
#    comments are appended/merged into existing comment cells, two new review
#    rows (15 & 16) are appended with a date + reviewer + new comments, and a
#    few new "F" (Comment) remarks are added.
#  - The active sheet/tab switches from "Review" to "HSI review".
#  - Review sheet's row heights for rows 3 and 5 grow (wrapped text got
#    longer) and the two new rows get explicit heights.

$wb = $excel.ActiveWorkbook
$review = $wb.Worksheets.Item("Review")
$hsi = $wb.Worksheets.Item("HSI review")

# ---------------------------------------------------------------------------
# Review sheet - point status / acceptance / comment updates
# ---------------------------------------------------------------------------

# Row 3 - point reviewed & closed, comment merged with the new follow-up note
$review.Range("E3").Value = "Closed"
$review.Range("F3").Value = @'
Mali 25/1/2020: 
Please remove "1.Project Name" no need for it
Mali6/2/2020: Point is reviewed and closed
'@

# Row 5 - new follow-up appended to the existing "Microcontroller block" note
$review.Range("F5").Value = @'
Mali 30/1/2020: Please provide a block "Microcontroller" has input side(Input switches) and output side (LEDs) 
Mali 6/2/2020: TI switch still a 3 phase switch while TI has 2 separate switches
'@

# Row 7 - accepted, closed, and a brand new comment
$review.Range("D7").Value = "Accepted"
$review.Range("E7").Value = "Closed"
$review.Range("F7").Value = "Please remove curly brackets"
$review.Range("F7").HorizontalAlignment = -4131
$review.Range("F7").VerticalAlignment = -4160

# Row 12 - new comment (still open)
$review.Range("F12").Value = "Mali 6/2/2020: Pointstill open"
$review.Range("F12").HorizontalAlignment = -4131
$review.Range("F12").VerticalAlignment = -4160

# Row 14 - accepted, closed, and a brand new comment
$review.Range("D14").Value = "Accepted"
$review.Range("E14").Value = "Closed"
$review.Range("F14").Value = "Mali 6/2/2020: Point is reviewed and closed"
$review.Range("F14").HorizontalAlignment = -4131
$review.Range("F14").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Review sheet - two brand new review rows (15 & 16)
# ---------------------------------------------------------------------------

$review.Range("A15").Value = 43984
$review.Range("A15").NumberFormat = "mm-dd-yy"
$review.Range("A15").HorizontalAlignment = -4108
$review.Range("A15").VerticalAlignment = -4108
$review.Range("B15").Value = "Ali"
$review.Range("B15").HorizontalAlignment = -4108
$review.Range("B15").VerticalAlignment = -4108
$review.Range("C15").Value = @'
Requirement "Req _ SEQUENCE OF TI FUNCTION RIGHT_PO5_CYRS_005-V01" & "Req _ SEQUENCE OF TI FUNCTION RIGHT_PO5_CYRS_006-V01", they don't declare that scenario of TI animation is repeated till TI switch is released.
'@
$review.Range("C15").HorizontalAlignment = -4131
$review.Range("C15").VerticalAlignment = -4160
$review.Range("C15").WrapText = $true
$review.Range("D15").HorizontalAlignment = -4108
$review.Range("D15").VerticalAlignment = -4108
$review.Range("D15").WrapText = $true
$review.Range("E15").Value = "Open"
$review.Range("E15").HorizontalAlignment = -4108
$review.Range("E15").VerticalAlignment = -4108
$review.Range("F15").HorizontalAlignment = -4131
$review.Range("F15").VerticalAlignment = -4160
$review.Rows.Item(15).RowHeight = 60

$review.Range("A16").Value = 43984
$review.Range("A16").NumberFormat = "mm-dd-yy"
$review.Range("A16").HorizontalAlignment = -4108
$review.Range("A16").VerticalAlignment = -4108
$review.Range("B16").Value = "Ali"
$review.Range("B16").HorizontalAlignment = -4108
$review.Range("B16").VerticalAlignment = -4108
$review.Range("C16").Value = @'
Req _ SELECT WELCOME MODE_PO5_CYRS_001-V01 doesn't declare
the switch is Mode switch it's just saying "If the switch"
'@
$review.Range("C16").HorizontalAlignment = -4131
$review.Range("C16").VerticalAlignment = -4160
$review.Range("C16").WrapText = $true
$review.Range("E16").Value = "Open"
$review.Range("E16").HorizontalAlignment = -4108
$review.Range("E16").VerticalAlignment = -4108
$review.Rows.Item(16).RowHeight = 45

# Extend the "Open/Closed" and "Accepted/Rejected" dropdown validations to
# cover the new rows.
$review.Range("E2:E16").Validation.Delete()
$review.Range("E2:E16").Validation.Add(3, 1, 1, "Open, Closed")
$review.Range("E2:E16").Validation.InputMessage = ""
$review.Range("E2:E16").Validation.ErrorMessage = ""

$review.Range("D2:D15").Validation.Delete()
$review.Range("D2:D15").Validation.Add(3, 1, 1, "Accepted, Rejected")
$review.Range("D2:D15").Validation.InputMessage = ""
$review.Range("D2:D15").Validation.ErrorMessage = ""

# Row heights that grew because the wrapped comment text got longer.
$review.Rows.Item(3).RowHeight = 45
$review.Rows.Item(5).RowHeight = 75

# ---------------------------------------------------------------------------
# Switch the active tab from "Review" to "HSI review"
# ---------------------------------------------------------------------------

$hsi.Activate()
